$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 74 data (215. Kth Largest Element in an Array / Heaps) ---
$ws.Range("A74").Value = "215. Kth Largest Element in an Array"
$ws.Range("B74").Value = "Medium"
$ws.Range("C74").Value = "Heaps"
$ws.Range("D74").Value = "The naive way is to sort and return len-k. We can solve with n + klogn time with a heap, but the optimal is average O(n) time with Quick Select. We choose a element as a pivot, and iterate the other elements, choosing to place it to the left or the right of the pivot based on the element we chose. Each side of the pivot is the `"partition`". At the end, we swap the pivot value with the element remaining at the pivot index. We recursively perform this on the partitions until we find the element at length-k."
$ws.Range("E74").Value = "https://leetcode.com/problems/kth-largest-element-in-an-array/solutions/60294/solution-explained/ "

# --- Copy formatting from the previous last row (row 73) so that fills/styles match ---
$ws.Range("A73").Copy()
$ws.Range("A74").PasteSpecial(-4122)
$ws.Range("B73").Copy()
$ws.Range("B74").PasteSpecial(-4122)
$ws.Range("C73").Copy()
$ws.Range("C74").PasteSpecial(-4122)
$ws.Range("D73").Copy()
$ws.Range("D74").PasteSpecial(-4122)
$ws.Range("E73").Copy()
$ws.Range("E74").PasteSpecial(-4122)

# --- Turn the new E74 cell into a live hyperlink pointing at the LeetCode solution ---
$ws.Hyperlinks.Add($ws.Range("E74"), "https://leetcode.com/problems/kth-largest-element-in-an-array/solutions/60294/solution-explained/ ")

# Re-apply the hyperlink-style formatting (Hyperlinks.Add resets the cell style)
$ws.Range("E73").Copy()
$ws.Range("E74").PasteSpecial(-4122)

# --- Grow table (Table2 / ListObject) to include the new row ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E74"))

# --- Update the sheet view to reflect the new scroll position / selection ---
$aw = $excel.ActiveWindow
$aw.ScrollRow = 49
$aw.ScrollColumn = 2
$ws.Range("E77").Select()
